$wb = $excel.ActiveWorkbook

# --- Sheet 1: Recommandations ---
$ws1 = $wb.Worksheets.Item("Recommandations")

$ws1.Cells.Item(2,1).Value = "BRVM - CONSOMMATION DE BASE     (**)"
$ws1.Cells.Item(2,2).Value = 0
$ws1.Cells.Item(2,3).Value = 3
$ws1.Cells.Item(2,4).Value = 761.47
$ws1.Cells.Item(2,5).Value = 261.68
$ws1.Cells.Item(2,6).Value = "🟡 Observer"
$ws1.Cells.Item(2,7).Value = "➖ Neutre"

$ws1.Cells.Item(3,1).Value = "BRVM-PRINCIPAL     (**)"
$ws1.Cells.Item(3,2).Value = 0
$ws1.Cells.Item(3,3).Value = 3
$ws1.Cells.Item(3,4).Value = 748.52
$ws1.Cells.Item(3,5).Value = 254.67
$ws1.Cells.Item(3,6).Value = "🟡 Observer"
$ws1.Cells.Item(3,7).Value = "➖ Neutre"

$ws1.Cells.Item(4,1).Value = "BRVM - INDUSTRIELS"
$ws1.Cells.Item(4,2).Value = 0
$ws1.Cells.Item(4,3).Value = 3
$ws1.Cells.Item(4,4).Value = 577.71
$ws1.Cells.Item(4,5).Value = 199.49
$ws1.Cells.Item(4,6).Value = "🟡 Observer"
$ws1.Cells.Item(4,7).Value = "➖ Neutre"

$ws1.Cells.Item(5,1).Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws1.Cells.Item(5,2).Value = 0
$ws1.Cells.Item(5,3).Value = 3
$ws1.Cells.Item(5,4).Value = 574.78
$ws1.Cells.Item(5,5).Value = 192.43
$ws1.Cells.Item(5,6).Value = "🟡 Observer"
$ws1.Cells.Item(5,7).Value = "➖ Neutre"

$ws1.Cells.Item(6,1).Value = "BRVM - SERVICES FINANCIERS"
$ws1.Cells.Item(6,2).Value = 0
$ws1.Cells.Item(6,3).Value = 3
$ws1.Cells.Item(6,4).Value = 475.56
$ws1.Cells.Item(6,5).Value = 159.58
$ws1.Cells.Item(6,6).Value = "🟡 Observer"
$ws1.Cells.Item(6,7).Value = "➖ Neutre"

$ws1.Cells.Item(7,1).Value = "BRVM-PRESTIGE"
$ws1.Cells.Item(7,2).Value = 0
$ws1.Cells.Item(7,3).Value = 3
$ws1.Cells.Item(7,4).Value = 460.38
$ws1.Cells.Item(7,5).Value = 153.94
$ws1.Cells.Item(7,6).Value = "🟡 Observer"
$ws1.Cells.Item(7,7).Value = "➖ Neutre"

$ws1.Cells.Item(8,1).Value = "BRVM – COMPOSITE TOTAL RETURN     (**)"
$ws1.Cells.Item(8,2).Value = 0
$ws1.Cells.Item(8,3).Value = 3
$ws1.Cells.Item(8,4).Value = 436.46
$ws1.Cells.Item(8,5).Value = 147.05
$ws1.Cells.Item(8,6).Value = "🟡 Observer"
$ws1.Cells.Item(8,7).Value = "➖ Neutre"

$ws1.Cells.Item(9,1).Value = "BRVM - ENERGIE"
$ws1.Cells.Item(9,2).Value = 0
$ws1.Cells.Item(9,3).Value = 3
$ws1.Cells.Item(9,4).Value = 372.58
$ws1.Cells.Item(9,5).Value = 125.81
$ws1.Cells.Item(9,6).Value = "🟡 Observer"
$ws1.Cells.Item(9,7).Value = "➖ Neutre"

$ws1.Cells.Item(10,1).Value = "BRVM - SERVICES PUBLICS"
$ws1.Cells.Item(10,2).Value = 0
$ws1.Cells.Item(10,3).Value = 3
$ws1.Cells.Item(10,4).Value = 356.9
$ws1.Cells.Item(10,5).Value = 119.94
$ws1.Cells.Item(10,6).Value = "🟡 Observer"
$ws1.Cells.Item(10,7).Value = "➖ Neutre"

$ws1.Cells.Item(11,1).Value = "BRVM - TELECOMMUNICATIONS"
$ws1.Cells.Item(11,2).Value = 0
$ws1.Cells.Item(11,3).Value = 3
$ws1.Cells.Item(11,4).Value = 298.85
$ws1.Cells.Item(11,5).Value = 100.2
$ws1.Cells.Item(11,6).Value = "🟡 Observer"
$ws1.Cells.Item(11,7).Value = "➖ Neutre"

$ws1.Cells.Item(12,1).Value = "EVIOSYS PACKAGING SIEM CI (SEMC)"
$ws1.Cells.Item(12,2).Value = 3
$ws1.Cells.Item(12,3).Value = 0
$ws1.Cells.Item(12,4).Value = 22.04
$ws1.Cells.Item(12,5).Value = 7.31
$ws1.Cells.Item(12,6).Value = "🟢 Achat"
$ws1.Cells.Item(12,7).Value = "✅ Renforcer"

$ws1.Cells.Item(13,1).Value = "UNIWAX CI (UNXC)"
$ws1.Cells.Item(13,2).Value = 3
$ws1.Cells.Item(13,3).Value = 0
$ws1.Cells.Item(13,4).Value = 21.95
$ws1.Cells.Item(13,5).Value = 7.29
$ws1.Cells.Item(13,6).Value = "🟢 Achat"
$ws1.Cells.Item(13,7).Value = "✅ Renforcer"

$ws1.Cells.Item(14,1).Value = "ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)"
$ws1.Cells.Item(14,2).Value = 3
$ws1.Cells.Item(14,3).Value = 0
$ws1.Cells.Item(14,4).Value = 21.84
$ws1.Cells.Item(14,5).Value = 6.94
$ws1.Cells.Item(14,6).Value = "🟢 Achat"
$ws1.Cells.Item(14,7).Value = "✅ Renforcer"

$ws1.Cells.Item(15,1).Value = "SETAO CI (STAC)"
$ws1.Cells.Item(15,2).Value = 3
$ws1.Cells.Item(15,3).Value = 0
$ws1.Cells.Item(15,4).Value = 21.81
$ws1.Cells.Item(15,5).Value = 7.33
$ws1.Cells.Item(15,6).Value = "🟢 Achat"
$ws1.Cells.Item(15,7).Value = "✅ Renforcer"

$ws1.Cells.Item(16,1).Value = "UNILEVER CI (UNLC)"
$ws1.Cells.Item(16,2).Value = 2
$ws1.Cells.Item(16,3).Value = 0
$ws1.Cells.Item(16,4).Value = 14.99
$ws1.Cells.Item(16,5).Value = 7.49
$ws1.Cells.Item(16,6).Value = "🟡 Observer"
$ws1.Cells.Item(16,7).Value = "➖ Neutre"

$ws1.Cells.Item(17,1).Value = "ORAGROUP TOGO (ORGT)"
$ws1.Cells.Item(17,2).Value = 1
$ws1.Cells.Item(17,3).Value = 1
$ws1.Cells.Item(17,4).Value = 4.84
$ws1.Cells.Item(17,5).Value = 7.45
$ws1.Cells.Item(17,6).Value = "🟡 Observer"
$ws1.Cells.Item(17,7).Value = "👀 À surveiller"

$ws1.Cells.Item(18,1).Value = "BANK OF AFRICA SENEGAL (BOAS)"
$ws1.Cells.Item(18,2).Value = 0
$ws1.Cells.Item(18,3).Value = 1
$ws1.Cells.Item(18,4).Value = -1.59
$ws1.Cells.Item(18,5).Value = -1.59
$ws1.Cells.Item(18,6).Value = "🟡 Observer"
$ws1.Cells.Item(18,7).Value = "➖ Neutre"

$ws1.Cells.Item(19,1).Value = "TOTALENERGIES MARKETING CI (TTLC)"
$ws1.Cells.Item(19,2).Value = 0
$ws1.Cells.Item(19,3).Value = 1
$ws1.Cells.Item(19,4).Value = -1.64
$ws1.Cells.Item(19,5).Value = -1.64
$ws1.Cells.Item(19,6).Value = "🟡 Observer"
$ws1.Cells.Item(19,7).Value = "➖ Neutre"

$ws1.Cells.Item(20,1).Value = "CORIS BANK INTERNATIONAL (CBIBF)"
$ws1.Cells.Item(20,2).Value = 0
$ws1.Cells.Item(20,3).Value = 1
$ws1.Cells.Item(20,4).Value = -2
$ws1.Cells.Item(20,5).Value = -2
$ws1.Cells.Item(20,6).Value = "🟡 Observer"
$ws1.Cells.Item(20,7).Value = "➖ Neutre"

$ws1.Cells.Item(21,1).Value = "AFRICA GLOBAL LOGISTICS CI (SDSC)"
$ws1.Cells.Item(21,2).Value = 0
$ws1.Cells.Item(21,3).Value = 1
$ws1.Cells.Item(21,4).Value = -2.29
$ws1.Cells.Item(21,5).Value = -2.29
$ws1.Cells.Item(21,6).Value = "🟡 Observer"
$ws1.Cells.Item(21,7).Value = "➖ Neutre"

$ws1.Cells.Item(22,1).Value = "VIVO ENERGY CI (SHEC)"
$ws1.Cells.Item(22,2).Value = 0
$ws1.Cells.Item(22,3).Value = 1
$ws1.Cells.Item(22,4).Value = -2.78
$ws1.Cells.Item(22,5).Value = -2.78
$ws1.Cells.Item(22,6).Value = "🟡 Observer"
$ws1.Cells.Item(22,7).Value = "➖ Neutre"

$ws1.Cells.Item(23,1).Value = "BERNABE CI (BNBC)"
$ws1.Cells.Item(23,2).Value = 0
$ws1.Cells.Item(23,3).Value = 1
$ws1.Cells.Item(23,4).Value = -2.94
$ws1.Cells.Item(23,5).Value = -2.94
$ws1.Cells.Item(23,6).Value = "🟡 Observer"
$ws1.Cells.Item(23,7).Value = "➖ Neutre"

$ws1.Cells.Item(24,1).Value = "SOGB CI (SOGC)"
$ws1.Cells.Item(24,2).Value = 0
$ws1.Cells.Item(24,3).Value = 1
$ws1.Cells.Item(24,4).Value = -4.34
$ws1.Cells.Item(24,5).Value = -4.34
$ws1.Cells.Item(24,6).Value = "🟡 Observer"
$ws1.Cells.Item(24,7).Value = "➖ Neutre"

$ws1.Cells.Item(25,1).Value = "BANK OF AFRICA BF (BOABF)"
$ws1.Cells.Item(25,2).Value = 0
$ws1.Cells.Item(25,3).Value = 2
$ws1.Cells.Item(25,4).Value = -4.62
$ws1.Cells.Item(25,5).Value = -2.72
$ws1.Cells.Item(25,6).Value = "🟡 Observer"
$ws1.Cells.Item(25,7).Value = "➖ Neutre"

$ws1.Cells.Item(26,1).Value = "ECOBANK COTE D''IVOIRE (ECOC)"
$ws1.Cells.Item(26,2).Value = 0
$ws1.Cells.Item(26,3).Value = 2
$ws1.Cells.Item(26,4).Value = -5.77
$ws1.Cells.Item(26,5).Value = -2.91
$ws1.Cells.Item(26,6).Value = "🟡 Observer"
$ws1.Cells.Item(26,7).Value = "➖ Neutre"

$ws1.Cells.Item(27,1).Value = "SODE CI (SDCC)"
$ws1.Cells.Item(27,2).Value = 0
$ws1.Cells.Item(27,3).Value = 1
$ws1.Cells.Item(27,4).Value = -6.89
$ws1.Cells.Item(27,5).Value = -6.89
$ws1.Cells.Item(27,6).Value = "🟡 Observer"
$ws1.Cells.Item(27,7).Value = "➖ Neutre"

$ws1.Cells.Item(28,1).Value = "FILTISAC CI (FTSC)"
$ws1.Cells.Item(28,2).Value = 0
$ws1.Cells.Item(28,3).Value = 1
$ws1.Cells.Item(28,4).Value = -7.26
$ws1.Cells.Item(28,5).Value = -7.26
$ws1.Cells.Item(28,6).Value = "🟡 Observer"
$ws1.Cells.Item(28,7).Value = "➖ Neutre"

$ws1.Cells.Item(29,1).Value = "NEI-CEDA CI (NEIC)"
$ws1.Cells.Item(29,2).Value = 0
$ws1.Cells.Item(29,3).Value = 1
$ws1.Cells.Item(29,4).Value = -7.49
$ws1.Cells.Item(29,5).Value = -7.49
$ws1.Cells.Item(29,6).Value = "🟡 Observer"
$ws1.Cells.Item(29,7).Value = "➖ Neutre"

# Remove now-obsolete trailing rows (30-33) so dimension shrinks to A1:G29
$ws1.Range("A30:G33").ClearContents()

# --- Sheet 2: Top_YTD ---
$ws2 = $wb.Worksheets.Item("Top_YTD")

$ws2.Cells.Item(2,1).Value = "BRVM - CONSOMMATION DE BASE     (**)"
$ws2.Cells.Item(2,2).Value = 4326.46

$ws2.Cells.Item(3,1).Value = "BRVM-PRINCIPAL     (**)"
$ws2.Cells.Item(3,2).Value = 4168.22

$ws2.Cells.Item(4,1).Value = "BRVM - INDUSTRIELS"
$ws2.Cells.Item(4,2).Value = 2401.4

$ws2.Cells.Item(5,1).Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws2.Cells.Item(5,2).Value = 2379.1

$ws2.Cells.Item(6,1).Value = "BRVM - SERVICES FINANCIERS"
$ws2.Cells.Item(6,2).Value = 1627.71

$ws2.Cells.Item(7,1).Value = "BRVM-PRESTIGE"
$ws2.Cells.Item(7,2).Value = 1528.26

$ws2.Cells.Item(8,1).Value = "BRVM – COMPOSITE TOTAL RETURN     (**)"
$ws2.Cells.Item(8,2).Value = 1379.31

$ws2.Cells.Item(9,1).Value = "BRVM - ENERGIE"
$ws2.Cells.Item(9,2).Value = 1026.8

$ws2.Cells.Item(10,1).Value = "BRVM - SERVICES PUBLICS"
$ws2.Cells.Item(10,2).Value = 949.8200000000001

$ws2.Cells.Item(11,1).Value = "BRVM - TELECOMMUNICATIONS"
$ws2.Cells.Item(11,2).Value = 695.4
